$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.766.83'
$ws.Range('E2').Value = '  -0.60%  '
$ws.Range('D3').Value = '1.895.17'
$ws.Range('E3').Value = '  +0.11%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').Value = '0.7646'
$ws.Range('E5').Value = '  +4.15%  '
$ws.Range('D6').Value = '240.12'
$ws.Range('E6').Value = '  -1.05%  '
$ws.Range('E7').Value = '  +0.08%  '
$ws.Range('D8').Value = '0.3034'
$ws.Range('E8').Value = '  -2.04%  '
$ws.Range('D9').Value = '25.36'
$ws.Range('E9').Value = '  -3.55%  '
$ws.Range('D10').Value = '0.06808'
$ws.Range('E10').Value = '  -1.29%  '
$ws.Range('D11').Value = '0.07968'
$ws.Range('E11').Value = '  +0.18%  '
$ws.Range('D12').Value = '1.895.82'
$ws.Range('E12').Value = '  +0.00%  '
$ws.Range('D13').Value = '0.7341'
$ws.Range('E13').Value = '  -4.89%  '
$ws.Range('D14').Value = '5.150'
$ws.Range('E14').Value = '  -1.36%  '
$ws.Range('D15').Value = '90.69'
$ws.Range('E15').Value = '  -0.91%  '
$ws.Range('D16').Value = '29.766.36'
$ws.Range('E16').Value = '  -0.65%  '
$ws.Range('D17').Value = '13.77'
$ws.Range('E17').Value = '  -2.57%  '
$ws.Range('E18').Value = '  +1.79%  '
$ws.Range('D19').Value = '240.79'
$ws.Range('E19').Value = '  +0.27%  '
$ws.Range('D21').Value = '1.0000'
$ws.Range('E21').Value = '  +0.03%  '
$ws.Range('B22').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C22').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D22').Value = '2.136.86'
$ws.Range('E22').Value = '  +0.88%  '
$ws.Range('B23').Value = 'BinanceUSD'
$ws.Range('C23').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D23').Value = '1.002'
$ws.Range('E23').Value = '  +0.13%  '
$ws.Range('B24').Value = 'Chainlink'
$ws.Range('C24').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D24').Value = '6.881'
$ws.Range('E24').Value = '  -0.74%  '
$ws.Range('B25').Value = 'Monero'
$ws.Range('C25').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D25').Value = '166.36'
$ws.Range('E25').Value = '  +0.85%  '
$ws.Range('B26').Value = 'Cosmos'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D26').Value = '9.197'
$ws.Range('E26').Value = '  -1.12%  '
$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').Value = '18.57'
$ws.Range('E27').Value = '  -1.41%  '
$ws.Range('B28').Value = 'Stellar'
$ws.Range('C28').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D28').Value = '0.1292'
$ws.Range('E28').Value = '  +1.63%  '
$ws.Range('B29').Value = 'LidoDAOToken'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D29').Value = '2.019'
$ws.Range('E29').Value = '  +0.08%  '
$ws.Range('B30').Value = 'Toncoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D30').Value = '1.402'
$ws.Range('E30').Value = '  +2.85%  '
$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D31').Value = '1.513'
$ws.Range('E31').Value = '  -1.22%  '
$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D32').Value = '4.254'
$ws.Range('E32').Value = '  -1.22%  '
$ws.Range('B33').Value = 'InternetComputer(DFINITY)'
$ws.Range('C33').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D33').Value = '4.055'
$ws.Range('E33').Value = '  -0.15%  '
$ws.Range('B34').Value = 'Hedera'
$ws.Range('C34').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D34').Value = '0.05190'
$ws.Range('E34').Value = '  +1.55%  '
$ws.Range('B35').Value = 'ARBITRUM'
$ws.Range('C35').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D35').Value = '1.244'
$ws.Range('E35').Value = '  -2.81%  '
$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D36').Value = '0.7227'
$ws.Range('E36').Value = '  -1.84%  '
$ws.Range('B37').Value = 'HuobiToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D37').Value = '2.717'
$ws.Range('E37').Value = '  +0.02%  '
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').Value = '0.01913'
$ws.Range('E38').Value = '  -0.38%  '
$ws.Range('B39').Value = 'MXToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D39').Value = '2.774'
$ws.Range('E39').Value = '  -0.23%  '
$ws.Range('B40').Value = 'FraxShare'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D40').Value = '6.135'
$ws.Range('E40').Value = '  -2.92%  '
$ws.Range('B41').Value = 'TheSandbox'
$ws.Range('C41').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D41').Value = '0.4388'
$ws.Range('E41').Value = '  -1.57%  '
$ws.Range('B42').Value = 'Aave'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D42').Value = '71.66'
$ws.Range('E42').Value = '  -3.45%  '
$ws.Range('B43').Value = 'PaxDollar'
$ws.Range('C43').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D43').Value = '1.001'
$ws.Range('E43').Value = '  +0.05%  '
$ws.Range('B45').Value = 'TrustWalletToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D45').Value = '0.8281'
$ws.Range('E45').Value = '  -1.31%  '
$ws.Range('B46').Value = 'Aptos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D46').Value = '7.591'
$ws.Range('E46').Value = '  -0.43%  '
$ws.Range('B47').Value = 'Quant'
$ws.Range('C47').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D47').Value = '99.62'
$ws.Range('E47').Value = '  -1.38%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').Value = '9.681'
$ws.Range('E48').Value = '  -1.07%  '
$ws.Range('B49').Value = 'RocketPoolETH'
$ws.Range('C49').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D49').Value = '2.037.95'
$ws.Range('E49').Value = '  -0.27%  '
$ws.Range('B50').Value = 'Elrond'
$ws.Range('C50').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D50').Value = '35.97'
$ws.Range('E50').Value = '  -1.69%  '
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').Value = '0.05922'
$ws.Range('E51').Value = '  -0.44%  '
